# Apply updated cryptocurrency price/volume data to the sheet.
# Column D ("Price") values are forced to remain plain text (matching the
# original inlineStr cell type) by prefixing a literal apostrophe, which is
# Excel's standard 'treat this as text' quote-prefix convention; this stops
# numeric-looking strings (e.g. "0.671") from being auto-converted to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "'43.951.43"
$ws.Cells.Item(2, 5).Value = "  +0.45%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "'2.354.48"
$ws.Cells.Item(3, 5).Value = "  +1.04%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.04%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'0.671"
$ws.Cells.Item(5, 5).Value = "  +4.16%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'235.74"
$ws.Cells.Item(6, 5).Value = "  +1.26%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'72.61"
$ws.Cells.Item(7, 5).Value = "  +10.24%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.08%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.563"
$ws.Cells.Item(9, 5).Value = "  +26.33%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.0987"
$ws.Cells.Item(10, 5).Value = "  +2.16%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "'27.88"
$ws.Cells.Item(11, 5).Value = "  +2.94%  "

# Row 12
$ws.Cells.Item(12, 2).Value = "TRON"
$ws.Cells.Item(12, 3).Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Cells.Item(12, 4).Value = "'0.107"
$ws.Cells.Item(12, 5).Value = "  +2.01%  "

# Row 13
$ws.Cells.Item(13, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(13, 4).Value = "'2.700.67"
$ws.Cells.Item(13, 5).Value = "  +1.23%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'16.85"
$ws.Cells.Item(14, 5).Value = "  +9.99%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'0.882"
$ws.Cells.Item(16, 5).Value = "  +5.14%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "'2.350.08"
$ws.Cells.Item(17, 5).Value = "  +0.76%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'43.825.38"
$ws.Cells.Item(18, 5).Value = "  +0.42%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  +2.95%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'76.37"
$ws.Cells.Item(20, 5).Value = "  +3.72%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +1.99%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'253.91"
$ws.Cells.Item(22, 5).Value = "  +2.18%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  -0.05%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'3.79"
$ws.Cells.Item(24, 5).Value = "  +0.48%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +6.71%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'2.29"
$ws.Cells.Item(27, 5).Value = "  +0.71%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'22.47"
$ws.Cells.Item(28, 5).Value = "  +0.56%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "'172.37"
$ws.Cells.Item(29, 5).Value = "  -1.39%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +10.03%  "

# Row 31
$ws.Cells.Item(31, 5).Value = "  -0.30%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'0.133"
$ws.Cells.Item(32, 5).Value = "  +5.44%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +4.74%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "'0.0716"
$ws.Cells.Item(34, 5).Value = "  +4.76%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +4.37%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +2.10%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +1.73%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'6.40"
$ws.Cells.Item(38, 5).Value = "  -2.20%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  +8.64%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'19.09"
$ws.Cells.Item(40, 5).Value = "  +8.71%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "BinanceUSD"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Cells.Item(41, 4).Value = "'1.00"
$ws.Cells.Item(41, 5).Value = "  -0.08%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "FraxShare"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(42, 4).Value = "'8.91"
$ws.Cells.Item(42, 5).Value = "  -1.95%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "'1.16"
$ws.Cells.Item(43, 5).Value = "  +1.04%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'0.0973"
$ws.Cells.Item(44, 5).Value = "  +2.43%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "Algorand"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(45, 4).Value = "'0.184"
$ws.Cells.Item(45, 5).Value = "  +14.13%  "

# Row 46
$ws.Cells.Item(46, 2).Value = "TrustWalletToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(46, 4).Value = "'1.21"
$ws.Cells.Item(46, 5).Value = "  +2.16%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "FTXToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(47, 4).Value = "'4.43"
$ws.Cells.Item(47, 5).Value = "  +0.71%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).Value = "'97.68"
$ws.Cells.Item(48, 5).Value = "  -0.97%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'1.434.00"
$ws.Cells.Item(49, 5).Value = "  -0.43%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'2.30"
$ws.Cells.Item(50, 5).Value = "  +1.25%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "HuobiToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(51, 4).Value = "'2.78"
$ws.Cells.Item(51, 5).Value = "  +1.70%  "
